{"js": "// Update the worksheet date and each \"three-digit \u00f7 one-digit\" problem\n// to the next day's generated set, per the commit diff.\nconst replacements = [\n  [\"2024-11-03 Sunday\", \"2024-11-04 Monday\"],\n  [\"382\u00f72=\", \"578\u00f78=\"],\n  [\"529\u00f78=\", \"501\u00f73=\"],\n  [\"704\u00f72=\", \"276\u00f77=\"],\n  [\"424\u00f78=\", \"454\u00f77=\"],\n  [\"145\u00f76=\", \"221\u00f75=\"],\n  [\"295\u00f73=\", \"778\u00f74=\"],\n  [\"332\u00f79=\", \"100\u00f73=\"],\n  [\"381\u00f75=\", \"235\u00f77=\"],\n  [\"606\u00f77=\", \"397\u00f78=\"],\n  [\"660\u00f76=\", \"123\u00f79=\"],\n  [\"885\u00f74=\", \"892\u00f77=\"],\n  [\"594\u00f72=\", \"893\u00f78=\"],\n  [\"766\u00f76=\", \"120\u00f74=\"],\n  [\"949\u00f75=\", \"185\u00f73=\"],\n  [\"271\u00f77=\", \"685\u00f78=\"],\n  [\"729\u00f76=\", \"491\u00f72=\"],\n  [\"710\u00f73=\", \"589\u00f76=\"],\n  [\"227\u00f76=\", \"240\u00f74=\"],\n  [\"780\u00f72=\", \"773\u00f73=\"],\n  [\"773\u00f78=\", \"602\u00f75=\"],\n  [\"947\u00f73=\", \"733\u00f73=\"],\n  [\"820\u00f79=\", \"870\u00f78=\"],\n  [\"215\u00f72=\", \"428\u00f76=\"],\n  [\"313\u00f76=\", \"202\u00f78=\"],\n  [\"918\u00f79=\", \"745\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each \"three-digit \u00f7 one-digit\" problem\n# to the next day's generated set, per the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-11-03 Sunday\", \"2024-11-04 Monday\"),\n  @(\"382\u00f72=\", \"578\u00f78=\"),\n  @(\"529\u00f78=\", \"501\u00f73=\"),\n  @(\"704\u00f72=\", \"276\u00f77=\"),\n  @(\"424\u00f78=\", \"454\u00f77=\"),\n  @(\"145\u00f76=\", \"221\u00f75=\"),\n  @(\"295\u00f73=\", \"778\u00f74=\"),\n  @(\"332\u00f79=\", \"100\u00f73=\"),\n  @(\"381\u00f75=\", \"235\u00f77=\"),\n  @(\"606\u00f77=\", \"397\u00f78=\"),\n  @(\"660\u00f76=\", \"123\u00f79=\"),\n  @(\"885\u00f74=\", \"892\u00f77=\"),\n  @(\"594\u00f72=\", \"893\u00f78=\"),\n  @(\"766\u00f76=\", \"120\u00f74=\"),\n  @(\"949\u00f75=\", \"185\u00f73=\"),\n  @(\"271\u00f77=\", \"685\u00f78=\"),\n  @(\"729\u00f76=\", \"491\u00f72=\"),\n  @(\"710\u00f73=\", \"589\u00f76=\"),\n  @(\"227\u00f76=\", \"240\u00f74=\"),\n  @(\"780\u00f72=\", \"773\u00f73=\"),\n  @(\"773\u00f78=\", \"602\u00f75=\"),\n  @(\"947\u00f73=\", \"733\u00f73=\"),\n  @(\"820\u00f79=\", \"870\u00f78=\"),\n  @(\"215\u00f72=\", \"428\u00f76=\"),\n  @(\"313\u00f76=\", \"202\u00f78=\"),\n  @(\"918\u00f79=\", \"745\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
